$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.917663037776947
$ws.Range("B1").Value = 2.844558000564575
$ws.Range("C1").Value = 8.741476058959961
$ws.Range("D1").Value = 2.027060270309448
$ws.Range("E1").Value = 1.143010377883911
